# Update countries & provincias Spain
# Applies the data refresh described in the commit to the "Pais" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last refreshed" timestamp banner in A1
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 18:43"

# ---------------------------------------------------------------------------
# 2) Two pairs of rows swapped rank because one country's numbers overtook
#    its neighbour's. Update the country name (column A) together with its
#    row of figures so the table keeps reflecting the correct country per
#    row.
# ---------------------------------------------------------------------------
$ws.Range("A64").Value = "Moldavia"
$ws.Range("A65").Value = "Kenia"

$ws.Range("A89").Value = "Zambia"
$ws.Range("A90").Value = "Malasia"

$ws.Range("A133").Value = "Tunez"
$ws.Range("A134").Value = "Benin"

$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# ---------------------------------------------------------------------------
# 3) Refresh the numeric figures (Casos totales, Nuevos casos, Casos
#    activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for every
#    affected row. Columns: B=Casos totales C=Nuevos casos D=Casos activos
#    E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes
# ---------------------------------------------------------------------------

function Set-Row($row, $b, $c, $d, $e, $g, $h) {
    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $d) { $ws.Cells.Item($row, 4).Value = $d }
    if ($null -ne $e) { $ws.Cells.Item($row, 5).Value = $e }
    if ($null -ne $g) { $ws.Cells.Item($row, 7).Value = $g }
    if ($null -ne $h) { $ws.Cells.Item($row, 8).Value = $h }
}

# Row 4 - Estados Unidos
Set-Row 4   5539841 10052 2904440 2462639 156   172762

# Row 17 - Argentina
Set-Row 17  $null   $null 211702  71741   20    5657

# Row 20 - Italia
Set-Row 20  253915  479   203786  14733   4     35396

# Row 22 - Alemania
Set-Row 22  224706  228   $null   12516   $null $null

# Row 43 - Bielorrusia
Set-Row 43  69516   92    67072   1834    3     610

# Row 58 - Argelia
Set-Row 58  38583   450   27017   10196   10    1370

# Row 59 - Suiza
Set-Row 59  $null   $null 33300   2833    $null $null

# Row 64 - now Moldavia (was Kenia), updated figures
Set-Row 64  30183   278   21220   8067    1     896

# Row 65 - now Kenia (was Moldavia), unchanged Kenia figures
Set-Row 65  30120   271   16656   12990   2     474

# Row 74 - Chequia
Set-Row 74  19956   65    13769   5790    2     397

# Row 85 - Senegal
Set-Row 85  12162   130   7677    4232    2     253

# Row 89 - now Zambia (was Malasia), updated figures
Set-Row 89  9343    157   8412    671     $null 260

# Row 90 - now Malasia (was Zambia), unchanged Malasia figures
Set-Row 90  9200    25    8859    216     $null 125

# Row 99 - Luxemburgo
Set-Row 99  7458    19    $null   835     $null $null

# Row 100 - Albania
Set-Row 100 7380    120   3794    3358    3     228

# Row 101 - Grecia
Set-Row 101 7075    217   $null   3043    2     228

# Row 105 - Republica de Yibuti
Set-Row 105 5369    2     5202    108     $null $null

# Row 125 - Mozambique
Set-Row 125 2855    64    1163    1673    $null $null

# Row 133 - now Tunez (was Benin), updated figures
Set-Row 133 2107    84    1358    695     $null 54

# Row 134 - now Benin (was Tunez), unchanged Benin figures
Set-Row 134 2063    $null 1690    334     $null 39

# Row 144 - Jordania
Set-Row 144 1378    39    1236    131     $null $null

# Row 173 - Birmania
Set-Row 173 375     1     329     40      $null $null

# Row 213 - now Islas Malvinas (was Montserrat), updated figures
Set-Row 213 $null   $null 13      $null   $null 0

# Row 214 - now Montserrat (was Islas Malvinas), unchanged Montserrat figures
Set-Row 214 $null   $null 12      $null   $null 1
